$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 165.54
$ws.Range("I15").Value = 165.54
$ws.Range("K15").Value = 496.62
$ws.Range("M15").Value = -327.62

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 324.25
$ws.Range("I33").Value = 438.4
$ws.Range("J33").Value = 134
$ws.Range("K33").Value = 438.4
$ws.Range("L33").Value = 134
$ws.Range("M33").Value = -209.4
$ws.Range("N33").Value = -592

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 66669560
$ws.Range("I94").Value = 71430960
$ws.Range("K94").Value = 71430960
$ws.Range("M94").Value = -71430509

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 782.0345
$ws.Range("I107").Value = 742.875
$ws.Range("K107").Value = 742.875
$ws.Range("M107").Value = 1177.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2300
$ws.Range("J113").Value = 2100
$ws.Range("L113").Value = 2100
$ws.Range("N113").Value = -8608

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1529.2858
$ws.Range("I116").Value = 2052.5
$ws.Range("J116").Value = 1320
$ws.Range("K116").Value = 2052.5
$ws.Range("L116").Value = 1320
$ws.Range("M116").Value = 1389.5
$ws.Range("N116").Value = -8204

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 759.4545000000001
$ws.Range("J129").Value = 905.7143
$ws.Range("L129").Value = 2717.1429
$ws.Range("N129").Value = -12717.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8942.637000000001
$ws.Range("I132").Value = 9756.9
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 29270.7
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = -26740.7
$ws.Range("N132").Value = -7460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 21278852
$ws.Range("I137").Value = 1391.88
$ws.Range("J137").Value = 45457784
$ws.Range("K137").Value = 4175.64
$ws.Range("L137").Value = 136373352
$ws.Range("M137").Value = -1625.64
$ws.Range("N137").Value = -136378452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12351.818
$ws.Range("I32").Value = 8929.887000000001
$ws.Range("J32").Value = 39727.273
$ws.Range("K32").Value = 8929.887000000001
$ws.Range("L32").Value = 39727.273
$ws.Range("M32").Value = -8642.887000000001
$ws.Range("N32").Value = -40301.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1102.56
$ws.Range("I45").Value = 994.625
$ws.Range("J45").Value = 1294.4445
$ws.Range("K45").Value = 994.625
$ws.Range("L45").Value = 1294.4445
$ws.Range("M45").Value = -617.625
$ws.Range("N45").Value = -2048.4445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12716.25
$ws.Range("I74").Value = 1570.7142
$ws.Range("J74").Value = 28320
$ws.Range("K74").Value = 1570.7142
$ws.Range("L74").Value = 28320
$ws.Range("M74").Value = -696.7141999999999
$ws.Range("N74").Value = -30068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 12716.25
$ws.Range("I77").Value = 1570.7142
$ws.Range("J77").Value = 28320
$ws.Range("K77").Value = 7853.571
$ws.Range("L77").Value = 141600
$ws.Range("M77").Value = -3485.571
$ws.Range("N77").Value = -150336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1853.3334
$ws.Range("I122").Value = 1853.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5560.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3110.0002
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1575.75
$ws.Range("I105").Value = 1426.125
$ws.Range("K105").Value = 1426.125
$ws.Range("M105").Value = 320.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 791.6667
$ws.Range("I16").Value = 737.5
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 737.5
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -450.5
$ws.Range("N16").Value = -1474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5211570.5
$ws.Range("I31").Value = 1416.1538
$ws.Range("J31").Value = 8776413
$ws.Range("K31").Value = 1416.1538
$ws.Range("L31").Value = 8776413
$ws.Range("M31").Value = -1121.1538
$ws.Range("N31").Value = -8777003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5211570.5
$ws.Range("I34").Value = 1416.1538
$ws.Range("J34").Value = 8776413
$ws.Range("K34").Value = 1416.1538
$ws.Range("L34").Value = 8776413
$ws.Range("M34").Value = -1214.1538
$ws.Range("N34").Value = -8776817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 33000
$ws.Range("J106").Value = 33000
$ws.Range("L106").Value = 33000
$ws.Range("N106").Value = -35524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 791.6667
$ws.Range("I113").Value = 737.5
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 737.5
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 1432.5
$ws.Range("N113").Value = -5240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 30389.758
$ws.Range("I12").Value = 100.166664
$ws.Range("J12").Value = 47698.094
$ws.Range("K12").Value = 300.499992
$ws.Range("L12").Value = 143094.282
$ws.Range("M12").Value = -127.499992
$ws.Range("N12").Value = -143440.282

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 672.25
$ws.Range("I86").Value = 399
$ws.Range("K86").Value = 1197
$ws.Range("M86").Value = -11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 672.25
$ws.Range("I89").Value = 399
$ws.Range("K89").Value = 3591
$ws.Range("M89").Value = 2337

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 698.625
$ws.Range("J92").Value = 721.5
$ws.Range("L92").Value = 2164.5
$ws.Range("N92").Value = -4660.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 62788.062
$ws.Range("I107").Value = 125301.625
$ws.Range("J107").Value = 41950.207
$ws.Range("K107").Value = 375904.875
$ws.Range("L107").Value = 125850.621
$ws.Range("M107").Value = -373984.875
$ws.Range("N107").Value = -129690.621

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2550.242
$ws.Range("J131").Value = 1657.3729
$ws.Range("L131").Value = 4972.1187
$ws.Range("N131").Value = -15052.1187

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -378
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 29000
$ws.Range("J109").Value = 29000
$ws.Range("L109").Value = 29000
$ws.Range("N109").Value = -31080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1169.3846
$ws.Range("J113").Value = 1185.5
$ws.Range("L113").Value = 1185.5
$ws.Range("N113").Value = -5525.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3914.5557
$ws.Range("I122").Value = 4116.7393
$ws.Range("J122").Value = 2752
$ws.Range("K122").Value = 12350.2179
$ws.Range("L122").Value = 8256
$ws.Range("M122").Value = -9900.2179
$ws.Range("N122").Value = -13156

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2172.1875
$ws.Range("I126").Value = 1850.5555
$ws.Range("J126").Value = 2585.7144
$ws.Range("K126").Value = 5551.666499999999
$ws.Range("L126").Value = 7757.1432
$ws.Range("M126").Value = -3081.666499999999
$ws.Range("N126").Value = -12697.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1670.2354
$ws.Range("I61").Value = 1408.8182
$ws.Range("J61").Value = 2149.5
$ws.Range("K61").Value = 1408.8182
$ws.Range("L61").Value = 2149.5
$ws.Range("M61").Value = -1206.8182
$ws.Range("N61").Value = -2553.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1819.7297
$ws.Range("I68").Value = 1619.6296
$ws.Range("J68").Value = 2360
$ws.Range("K68").Value = 1619.6296
$ws.Range("L68").Value = 2360
$ws.Range("M68").Value = -870.6296
$ws.Range("N68").Value = -3858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1819.7297
$ws.Range("I71").Value = 1619.6296
$ws.Range("J71").Value = 2360
$ws.Range("K71").Value = 8098.148
$ws.Range("L71").Value = 11800
$ws.Range("M71").Value = -4354.148
$ws.Range("N71").Value = -19288

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1670.2354
$ws.Range("I113").Value = 1408.8182
$ws.Range("J113").Value = 2149.5
$ws.Range("K113").Value = 1408.8182
$ws.Range("L113").Value = 2149.5
$ws.Range("M113").Value = 761.1818000000001
$ws.Range("N113").Value = -6489.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 2765.261
$ws.Range("I22").Value = 550.5
$ws.Range("J22").Value = 2976.1904
$ws.Range("K22").Value = 550.5
$ws.Range("L22").Value = 2976.1904
$ws.Range("M22").Value = -257.5
$ws.Range("N22").Value = -3562.1904

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 34988.5
$ws.Range("J109").Value = 34988.5
$ws.Range("L109").Value = 34988.5
$ws.Range("N109").Value = -37762.5
